$d = $word.ActiveDocument

# Make sure we don't end up recording this whole rewrite as tracked
# insertions/deletions.
$d.TrackRevisions = $false

# Replace the entire body (all paragraphs, keeping the existing sectPr)
# with the new, clean content described by the commit: the heading is
# simplified, the old paragraph (full of tracked changes / comments) is
# replaced by three new plain paragraphs, and a final paragraph holding
# the _GoBack bookmark + trailing space is kept.
$newBodyXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Valget for repository patternet</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Et repository pattern er et abstraktionslag til databasen for at simplificere koden. Da patternet giver et ekstra abstraktionslag kan det benytte på alle databaser selv om de bruger Entity</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Framework eller ADO.NET, el. lignende, da BLL kalder metoder gennem interfaces. Derudover giver repository pattern mulighed for at unit teste BLL i stedet for at integrations teste det, netop på grund af det abstraktionslag som mockes ud.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Der blev valgt at implementere et repository pattern i pristjek220 for at kunne unit teste programmets BLL. </w:t></w:r><w:r><w:t>Repository patternets abstraktionslag giver også pristjek220 mulighed for at kunne benytte sig af databaser som ikke bruger Entity Framework, så hvis der i fremtiden skulle blive udarbejdet noget nyere og bedre eller der bare bliver besluttet at det ikke skal køre på Enity Frameworket længere så kan programmet nemt skiftes over.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@

$d.Content.InsertXML($newBodyXml)

# Strip any comments left over from the original draft (their anchors are
# gone now anyway since the paragraph was rewritten, but make sure the
# comments collection itself is empty too).
for ($i = $d.Comments.Count; $i -ge 1; $i--) {
    $d.Comments.Item($i).Delete()
}
